# Case and Fatality Demographics Data Updated
# Updates the raw case counts (column B) on each of the six demographic
# breakdown sheets. All of the percentage figures in column C, and the
# "Total"/"Grand Total" rows, are formulas and recalculate automatically.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Cases by Age Group
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Age Group")
$ws.Range("B3").Value  = 1382
$ws.Range("B4").Value  = 3821
$ws.Range("B5").Value  = 15806
$ws.Range("B6").Value  = 17337
$ws.Range("B7").Value  = 15208
$ws.Range("B8").Value  = 12831
$ws.Range("B9").Value  = 4645
$ws.Range("B10").Value = 3143
$ws.Range("B12").Value = 1254
$ws.Range("B13").Value = 1950
$ws.Range("E3").Select()

# ---------------------------------------------------------------------------
# Sheet: Cases by Gender
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Gender")
$ws.Range("B2").Value = 27126
$ws.Range("B3").Value = 51549
$ws.Range("B4").Select()

# ---------------------------------------------------------------------------
# Sheet: Cases by RaceEthnicity
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws.Range("B2").Value = 965
$ws.Range("B3").Value = 13062
$ws.Range("B4").Value = 28490
$ws.Range("B5").Value = 579
$ws.Range("B6").Value = 27788
$ws.Range("B7").Value = 8689
$ws.Range("B8").Select()

# ---------------------------------------------------------------------------
# Sheet: Fatalities by Age Group
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Age Group")
$ws.Range("B4").Value  = 34
$ws.Range("B5").Value  = 261
$ws.Range("B6").Value  = 868
$ws.Range("B7").Value  = 2528
$ws.Range("B8").Value  = 5735
$ws.Range("B9").Value  = 4775
$ws.Range("B10").Value = 6140
$ws.Range("B11").Value = 6768
$ws.Range("B12").Value = 6676
$ws.Range("B13").Value = 16785
$ws.Range("C8").Select()

# ---------------------------------------------------------------------------
# Sheet: Fatalities by Gender
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Gender")
$ws.Range("B2").Value = 21204
$ws.Range("B3").Value = 29385
$ws.Range("E16").Select()

# ---------------------------------------------------------------------------
# Sheet: Fatalities by Race-Ethnicity
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws.Range("B2").Value = 1080
$ws.Range("B3").Value = 5129
$ws.Range("B4").Value = 23500
$ws.Range("B5").Value = 274
$ws.Range("B6").Value = 20584
$ws.Range("D15").Select()

# Re-activate "Fatalities by Age Group" last so it ends up the active tab,
# matching the saved workbook's activeTab/tabSelected state.
$wb.Worksheets.Item("Fatalities by Age Group").Activate()
